$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column, new text value.
# The price/volume columns (D, E) hold numeric-looking text (e.g. "1.007", "28.139.66")
# that must stay plain text, so we force NumberFormat "@" before assigning the value and
# restore the original cell style afterwards to avoid leaving stray formatting behind.
$updates = @(
    ,@(2, 4, '28.139.66')
    ,@(2, 5, '  +2.67%  ')
    ,@(3, 4, '1.822.38')
    ,@(3, 5, '  +1.42%  ')
    ,@(4, 4, '1.007')
    ,@(4, 5, '  +0.21%  ')
    ,@(5, 4, '340.14')
    ,@(5, 5, '  +0.71%  ')
    ,@(6, 4, '1.003')
    ,@(6, 5, '  +0.21%  ')
    ,@(7, 4, '0.3936')
    ,@(7, 5, '  +3.68%  ')
    ,@(8, 4, '0.3510')
    ,@(8, 5, '  +1.51%  ')
    ,@(9, 4, '48.19')
    ,@(9, 5, '  -0.75%  ')
    ,@(10, 4, '1.203')
    ,@(10, 5, '  +0.14%  ')
    ,@(11, 4, '0.07614')
    ,@(11, 5, '  +1.70%  ')
    ,@(12, 4, '1.003')
    ,@(12, 5, '  +0.16%  ')
    ,@(13, 4, '22.25')
    ,@(13, 5, '  +0.63%  ')
    ,@(14, 4, '6.555')
    ,@(14, 5, '  +1.32%  ')
    ,@(15, 4, '1.826.82')
    ,@(15, 5, '  +1.76%  ')
    ,@(16, 4, '7.224')
    ,@(16, 5, '  +2.02%  ')
    ,@(17, 4, '0.00001110')
    ,@(17, 5, '  +0.90%  ')
    ,@(18, 4, '0.06747')
    ,@(18, 5, '  +1.33%  ')
    ,@(19, 4, '85.64')
    ,@(19, 5, '  +1.19%  ')
    ,@(20, 4, '1.002')
    ,@(20, 5, '  +0.10%  ')
    ,@(21, 4, '17.88')
    ,@(21, 5, '  +3.19%  ')
    ,@(22, 4, '6.636')
    ,@(22, 5, '  +1.80%  ')
    ,@(23, 4, '28.144.56')
    ,@(23, 5, '  +2.74%  ')
    ,@(24, 4, '12.69')
    ,@(24, 5, '  +1.23%  ')
    ,@(25, 4, '2.408')
    ,@(25, 5, '  -0.79%  ')
    ,@(26, 4, '2.566')
    ,@(26, 5, '  +0.31%  ')
    ,@(27, 4, '1.490')
    ,@(27, 5, '  -0.11%  ')
    ,@(28, 4, '21.50')
    ,@(28, 5, '  +0.53%  ')
    ,@(29, 4, '155.62')
    ,@(29, 5, '  +2.65%  ')
    ,@(30, 4, '2.033.35')
    ,@(30, 5, '  +1.73%  ')
    ,@(31, 4, '136.72')
    ,@(31, 5, '  +2.07%  ')
    ,@(32, 4, '6.205')
    ,@(32, 5, '  +1.53%  ')
    ,@(33, 4, '4.050')
    ,@(33, 5, '  -0.08%  ')
    ,@(34, 4, '0.08849')
    ,@(34, 5, '  +1.94%  ')
    ,@(35, 4, '13.17')
    ,@(35, 5, '  -0.46%  ')
    ,@(36, 4, '5.537')
    ,@(36, 5, '  +1.75%  ')
    ,@(37, 2, 'Hedera')
    ,@(37, 3, 'https://coinranking.com/coin/jad286TjB+hedera-hbar')
    ,@(37, 4, '0.06620')
    ,@(37, 5, '  +4.05%  ')
    ,@(38, 2, 'TheSandbox')
    ,@(38, 3, 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand')
    ,@(38, 4, '0.6997')
    ,@(38, 5, '  +1.39%  ')
    ,@(39, 2, 'VeChain')
    ,@(39, 3, 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet')
    ,@(39, 4, '0.02447')
    ,@(39, 5, '  +4.83%  ')
    ,@(40, 5, '  -2.05%  ')
    ,@(41, 4, '0.2237')
    ,@(41, 5, '  +1.66%  ')
    ,@(42, 4, '1.270')
    ,@(42, 5, '  -0.58%  ')
    ,@(43, 4, '8.551')
    ,@(43, 5, '  -3.40%  ')
    ,@(44, 4, '14.80')
    ,@(44, 5, '  +2.16%  ')
    ,@(45, 4, '0.6518')
    ,@(45, 5, '  +1.25%  ')
    ,@(46, 4, '3.891')
    ,@(46, 5, '  +0.66%  ')
    ,@(47, 4, '2.176')
    ,@(47, 5, '  +2.27%  ')
    ,@(48, 4, '131.93')
    ,@(48, 5, '  +1.10%  ')
    ,@(49, 4, '0.07232')
    ,@(49, 5, '  +0.38%  ')
    ,@(50, 4, '80.50')
    ,@(50, 5, '  +0.92%  ')
    ,@(51, 4, '1.255')
    ,@(51, 5, '  +2.84%  ')
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $text = $u[2]
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}
